# Regenerate save_data to use K (column G) instead of Strike#.
# For each data row (2-35) recompute/write the new K value into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 9
    3  = 5
    4  = 2
    5  = 7
    6  = 8
    7  = 6
    8  = 8
    9  = 8
    10 = 7
    11 = 4
    12 = 8
    13 = 6
    14 = 4
    15 = 7
    16 = 6
    17 = 3
    18 = 10
    19 = 7
    20 = 14
    21 = 13
    22 = 5
    23 = 7
    24 = 4
    25 = 4
    26 = 9
    27 = 9
    28 = 5
    29 = 5
    30 = 4
    31 = 4
    32 = 3
    33 = 7
    34 = 3
    35 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
